# feat: add 2022-Q1 data
#
# 1. Duplicate the "2021-Q3" sheet (same column layout/styling) right after it,
#    rename the copy to "2022-Q1", and fill in the new quarter's fund data.
# 2. Insert a new leading data row in "总计" for "2022-Q1" and renumber the
#    existing index column.

$wb = $excel.ActiveWorkbook

# --- 1. Build the new "2022-Q1" sheet from the "2021-Q3" template ---------
$q3 = $wb.Worksheets.Item("2021-Q3")
$q3.Copy($null, $q3)
$q1 = $wb.Worksheets.Item(3)
$q1.Name = "2022-Q1"

# Header row: "基金金额" -> "基金规模" (everything else stays the same)
$q1.Range("D1").Value = "基金规模"

# Row 2 data (fund-code / monetary-looking values are stored as text, same
# convention the source workbook already uses for this table)
$q1.Range("B2:G2").NumberFormat = "@"
$q1.Range("B2").Value = "968029"
$q1.Range("C2").Value = "恒生指数基金M类人民币（对冲）份额"
$q1.Range("D2").Value = "25.09"
$q1.Range("E2").Value = "97.94"
$q1.Range("F2").Value = "4.04"
$q1.Range("G2").Value = "1.0136"
$q1.Range("H2").Value = 7

# --- 2. Update the "总计" summary sheet ------------------------------------
$total = $wb.Worksheets.Item("总计")

# Push the existing data rows down (row 3 -> row 4, row 2 -> row 3) to make
# room for the new "2022-Q1" row, writing values directly so the
# un-styled B:D columns stay un-styled (matches the original layout).
$total.Range("B4").Value = "2021-Q2"
$total.Range("C4").Value = 3
$total.Range("D4").Value = 1.65

$total.Range("B3").Value = "2021-Q3"
$total.Range("C3").Value = 1
$total.Range("D3").Value = 1.16

$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 1
$total.Range("D2").Value = 1.01

$total.Range("A2").Value = 0
$total.Range("A3").Value = 1
$total.Range("A4").Value = 2

# Column A keeps its bold/bordered "index" style on every data row.
$total.Range("A2").Copy()
$total.Range("A3:A4").PasteSpecial(-4122)
